$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdFileName = "b683438f-5d41-41ae-8334-f1fd00e193ac.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7bbeb870c051d3247655ec8189e10b3f70044875/e2e/b683438f-5d41-41ae-8334-f1fd00e193ac.md"

# --- Overview sheet: Status columns for zh-cn (E2) and de-de (F2) ---
$ws1.Range("E2").Value2 = $statusText
$ws1.Range("F2").Value2 = $statusText

# --- zh-cn sheet: Status column (C2) mirrors the same text ---
$ws2.Range("C2").Value2 = $statusText

# --- de-de sheet: Status column (C2) mirrors the same text ---
$ws3.Range("C2").Value2 = $statusText

# --- zh-cn sheet: Latest Target File (I2), Latest Handback File (J2), Latest Handback DateTime (K2) ---
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdFileName)
$ws2.Range("I2").Style = "HyperLink"
$ws2.Range("J2").Value2 = "b683438f-5d41-41ae-8334-f1fd00e193ac.cce209a9e16e8c4ac9a0f710c55f14bb4dddf846.zh-cn.xlf"
$ws2.Range("K2").Value2 = "2016-08-31 09:14:32"

# --- de-de sheet: Latest Target File (I2), Latest Handback File (J2), Latest Handback DateTime (K2) ---
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdFileName)
$ws3.Range("I2").Style = "HyperLink"
$ws3.Range("J2").Value2 = "b683438f-5d41-41ae-8334-f1fd00e193ac.cce209a9e16e8c4ac9a0f710c55f14bb4dddf846.de-de.xlf"
$ws3.Range("K2").Value2 = "2016-08-31 09:14:40"

# --- Column width adjustments (widened to fit the longer handback strings) ---
$ws1.Columns.Item(5).ColumnWidth = 29.1   # Overview col E
$ws1.Columns.Item(6).ColumnWidth = 29.1   # Overview col F

$ws2.Columns.Item(3).ColumnWidth = 29.1   # zh-cn col C (Status)
$ws2.Columns.Item(9).ColumnWidth = 39.17  # zh-cn col I (Latest Target File)
$ws2.Columns.Item(10).ColumnWidth = 39.17 # zh-cn col J (Latest Handback File)

$ws3.Columns.Item(3).ColumnWidth = 29.1   # de-de col C (Status)
$ws3.Columns.Item(9).ColumnWidth = 39.17  # de-de col I (Latest Target File)
$ws3.Columns.Item(10).ColumnWidth = 39.17 # de-de col J (Latest Handback File)
